# Weekly data update: insert a new price record for "Kiwi" / "Hayward" /
# "Primera" quality (Región de O'Higgins origin) ahead of the existing
# history, pushing the existing rows 249-356 down to 250-357.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 249; this shifts rows 249:356 down to
# 250:357 (and the sheet dimension grows from T356 to T357 automatically).
$ws.Rows("249:249").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A249").Value = 7
$ws.Range("B249").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C249").Value = "Ñuble"
$ws.Range("D249").Value = 45141
$ws.Range("E249").Value = 16
$ws.Range("F249").Value = "Fruta"
$ws.Range("G249").Value = 100101
$ws.Range("H249").Value = "Berries"
$ws.Range("I249").Value = 100101007
$ws.Range("J249").Value = "Kiwi"
$ws.Range("K249").Value = "Hayward"
$ws.Range("L249").Value = "Primera"
$ws.Range("M249").Value = 40
$ws.Range("N249").Value = 14000
$ws.Range("O249").Value = 14000
$ws.Range("P249").Value = 14000
$ws.Range("Q249").Value = "$/bandeja 18 kilos"
$ws.Range("R249").Value = "Región de O'Higgins"
$ws.Range("S249").Value = 778
$ws.Range("T249").Value = 18
